$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "C"="2.177054233802296"; "D"="2.94415561582875"; "E"="16.66802848903047"; "F"="24.47134540445464"; "G"="3.567526063762902"; "I"="19.59752123846184"; "O"="20.89182271100765" }
    3 = @{ "C"="2.172435147857623"; "D"="2.954999009823079"; "E"="15.70626053483033"; "F"="23.77279097030631"; "G"="3.571234858655119"; "I"="19.17331663189021"; "O"="20.41879351793053" }
    4 = @{ "C"="2.169938876734689"; "D"="2.962158412889017"; "E"="15.09050342618857"; "F"="23.34195226825212"; "G"="3.573629020972074"; "I"="18.9147953252135"; "O"="20.12992364643339" }
    5 = @{ "C"="2.169007673702026"; "D"="2.965201582666258"; "E"="14.83349963441479"; "F"="23.16620146761176"; "G"="3.574634190340249"; "I"="18.81009776041389"; "O"="20.01278767864559" }
    6 = @{ "C"="2.168858261243664"; "D"="2.965714475023161"; "E"="14.79046609958024"; "F"="23.13701640717184"; "G"="3.574802884975814"; "I"="18.79275700330085"; "O"="19.99337790370625" }
    7 = @{ "C"="2.169925969035361"; "D"="2.962198946007249"; "E"="15.08706160069392"; "F"="23.33958235184757"; "G"="3.573642457318738"; "I"="18.91338048007337"; "O"="20.12834131072605" }
    8 = @{ "C"="2.175391400923914"; "D"="2.947790048412123"; "E"="16.34179060929392"; "F"="24.23106040309863"; "G"="3.568780657122227"; "I"="19.45094988731148"; "O"="20.72850622303068" }
    9 = @{ "C"="2.188778008507341"; "D"="2.923534966248326"; "E"="18.73688618351986"; "F"="25.95145186653721"; "G"="3.56016906549611"; "I"="20.51348828274725"; "O"="21.90996168125045" }
    10 = @{ "C"="2.200200076093151"; "D"="2.908184729423083"; "E"="20.41905478277341"; "F"="27.18294755802046"; "G"="3.554396683160728"; "I"="21.29013090498095"; "O"="22.77054280716389" }
    11 = @{ "C"="2.205731232051471"; "D"="2.901745097343441"; "E"="21.14206437904656"; "F"="27.73323559559246"; "G"="3.551889431349227"; "I"="21.64077404586611"; "O"="23.15842820242476" }
    12 = @{ "C"="2.207873030907799"; "D"="2.899385351457923"; "E"="21.40981913813112"; "F"="27.93997367935109"; "G"="3.550956931082514"; "I"="21.77303194831455"; "O"="23.30463903254853" }
    13 = @{ "C"="2.207409669198076"; "D"="2.899890048953038"; "E"="21.35242109623705"; "F"="27.89552529305414"; "G"="3.551157010091559"; "I"="21.74457326096305"; "O"="23.27318221303181" }
    14 = @{ "C"="2.20590649727083"; "D"="2.901549376282046"; "E"="21.16421342581507"; "F"="27.75027818211138"; "G"="3.551812375130762"; "I"="21.65166624362909"; "O"="23.17047139746499" }
    15 = @{ "C"="2.204991888575427"; "D"="2.902576045441826"; "E"="21.04814619265929"; "F"="27.66108993610327"; "G"="3.552216007939524"; "I"="21.59468585420252"; "O"="23.10746597431885" }
    16 = @{ "C"="2.199845268635152"; "D"="2.908616564002022"; "E"="20.37095809521562"; "F"="27.14676622563292"; "G"="3.554562915103452"; "I"="21.2671498507566"; "O"="22.7451076581375" }
    17 = @{ "C"="2.196773189785387"; "D"="2.912461845258358"; "E"="19.94473900356369"; "F"="26.82854251495648"; "G"="3.556032967714196"; "I"="21.06543382608179"; "O"="22.52177723928166" }
    18 = @{ "C"="2.19503782080948"; "D"="2.914724641524839"; "E"="19.69561312536637"; "F"="26.64458609276428"; "G"="3.556889676510139"; "I"="20.94917007532184"; "O"="22.39299357549156" }
    19 = @{ "C"="2.194455711243687"; "D"="2.915499540319687"; "E"="19.61057936871486"; "F"="26.58215009205562"; "G"="3.557181666024127"; "I"="20.90976804047983"; "O"="22.34933788662767" }
    20 = @{ "C"="2.197096952948382"; "D"="2.912047215823319"; "E"="19.99052197472227"; "F"="26.8625150876032"; "G"="3.555875322552377"; "I"="21.08693294040409"; "O"="22.54558644856849" }
    21 = @{ "C"="2.206346740066932"; "D"="2.901059847334213"; "E"="21.21965802419764"; "F"="27.79298702711993"; "G"="3.55161941968067"; "I"="21.67897056426082"; "O"="23.20065950705425" }
    22 = @{ "C"="2.212667072012585"; "D"="2.894338630682327"; "E"="21.9878201978673"; "F"="28.39143486394574"; "G"="3.548936635058832"; "I"="22.06279337296465"; "O"="23.62480026203088" }
    23 = @{ "C"="2.209268943687952"; "D"="2.897883577502841"; "E"="21.58104139471089"; "F"="28.07298321618889"; "G"="3.550359496084197"; "I"="21.85826868709069"; "O"="23.39884195881071" }
    24 = @{ "C"="2.196950483585782"; "D"="2.912234507584557"; "E"="19.96983622591272"; "F"="26.84715920235389"; "G"="3.555946557926888"; "I"="21.07721409737909"; "O"="22.5348235017563" }
    25 = @{ "C"="2.184874773439678"; "D"="2.929665206824606"; "E"="18.07972582065041"; "F"="25.49069301138038"; "G"="3.562400784297299"; "I"="20.22611050585472"; "O"="21.590943052209" }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = [double]$rowData[$col]
    }
}
